$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price (D) values are plain text in the source data (e.g. "0.570",
# "59.350.28" as thousands-grouped). Cells whose new text parses as a plain
# number need NumberFormat "@" set first so Excel stores them as Text -
# matching the original inlineStr cells - instead of silently coercing
# them into numbers (which would also eat significant trailing zeros).

# Row 2
$ws.Range("D2").Value = "59.349.26"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "2.639.50"
$ws.Range("E3").Value = "  -0.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.55"
$ws.Range("E5").Value = "  +1.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.02"
$ws.Range("E6").Value = "  -1.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  -3.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("E10").Value = "  +1.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("E12").Value = "  +1.03%  "

# Row 13
$ws.Range("D13").Value = "3.109.96"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("D14").Value = "59.319.38"
$ws.Range("E14").Value = "  +0.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.94"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("D17").Value = "2.614.89"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.97"
$ws.Range("E18").Value = "  +0.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.46"
$ws.Range("E19").Value = "  +0.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.62"
$ws.Range("E20").Value = "  +3.03%  "

# Row 21
$ws.Range("E21").Value = "  +1.58%  "

# Row 22
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.73"
$ws.Range("E23").Value = "  +3.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.418"
$ws.Range("E24").Value = "  +1.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +0.63%  "

# Row 26
$ws.Range("E26").Value = "  -0.38%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.26"
$ws.Range("E27").Value = "  +1.90%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0802"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("E29").Value = "  -3.14%  "

# Row 30
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("E31").Value = "  +1.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.07"
$ws.Range("E32").Value = "  +1.81%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.89"
$ws.Range("E33").Value = "  +0.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.20"
$ws.Range("E34").Value = "  +1.26%  "

# Row 35
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.891"
$ws.Range("E36").Value = "  -0.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.866"
$ws.Range("E37").Value = "  -1.56%  "

# Row 38
$ws.Range("E38").Value = "  +0.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.66"
$ws.Range("E39").Value = "  -0.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.66"
$ws.Range("E40").Value = "  +2.35%  "

# Row 41
$ws.Range("E41").Value = "  -0.31%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  -4.06%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0975"
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "272.25"
$ws.Range("E44").Value = "  -1.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.41"
$ws.Range("E45").Value = "  -1.75%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0539"
$ws.Range("E46").Value = "  +0.67%  "

# Row 48
$ws.Range("D48").Value = "2.039.12"
$ws.Range("E48").Value = "  -0.52%  "

# Row 49
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0230"
$ws.Range("E50").Value = "  +0.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.92"
$ws.Range("E51").Value = "  +0.00%  "
